# Atualização de bases das ligas, do dia: 28-05-2024 às 19:13
#
# The underlying source data re-sorted / re-emitted a handful of match rows on
# the "Portugal Segunda Liga" sheet. For a set of adjacent row pairs the two
# matches traded places (everything except the leading row-index column A
# swapped between the two rows). Swapping the data this way also causes the
# "FC Porto B" / "Academico Viseu" shared-string table entries to end up
# re-ordered relative to each other, which is why plain references to those
# two team names shift around the sheet - that part is just bookkeeping that
# Excel keeps consistent automatically whenever cell text is (re)written, so
# we only need to move the real match data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portugal Segunda Liga")

# (rowA, rowB) pairs whose B:AD contents (match id, teams, odds, results, ...)
# got swapped, keeping column A (the running index) fixed on each row.
$pairs = @(
    @(5, 6),
    @(32, 33),
    @(134, 135),
    @(140, 141),
    @(151, 152),
    @(186, 187),
    @(201, 202),
    @(221, 223),
    @(241, 242),
    @(260, 261),
    @(278, 279),
    @(293, 294),
    @(296, 297),
    @(303, 304),
    @(306, 307)
)

foreach ($pair in $pairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]

    $rangeA = $ws.Range("B$rowA`:AD$rowA")
    $rangeB = $ws.Range("B$rowB`:AD$rowB")

    $valuesA = $rangeA.Value2
    $valuesB = $rangeB.Value2

    $rangeA.Value2 = $valuesB
    $rangeB.Value2 = $valuesA
}
